# Resume update: add Docker skill (and Java / Git skills), and tidy up the
# "Present Day" run-split that was left over from a prior edit (the cursor
# bookmark is relocated to the last real edit, which is the Machine
# Learning/NLP line).
$d = $word.ActiveDocument

function Replace-ParagraphXml($findText, $bodyXml) {
    $rng = $d.Content
    $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $findText"
    }
    $para = $rng.Paragraphs(1).Range
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $bodyXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData>' +
        '</pkg:part>' +
        '</pkg:package>'
    $para.InsertXML($xml)
}

# 1) "- Programming Languages: Python, Ruby, C/C++"
#    -> split out "Java, " as its own run between Python and Ruby.
Replace-ParagraphXml " Python, Ruby, C/C++" (
    '<w:p><w:pPr><w:pStyle w:val="17"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr>' +
    '<w:spacing w:after="0" w:afterLines="0" w:line="240" w:lineRule="auto"/><w:ind w:left="360" w:leftChars="0"/>' +
    '<w:rPr><w:rFonts w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:b/><w:bCs w:val="0"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>- Programming Languages:</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> Python, </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Java, </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Ruby, C/C++</w:t></w:r>' +
    '</w:p>'
)

# 2) "- Artificial Intelligence: Machine Learning, Natural Language Processing"
#    -> split into "Machine Learning" / (bookmark _GoBack moves here) / ", Natural Language Processing"
Replace-ParagraphXml " Machine Learning, Natural Language Processing" (
    '<w:p><w:pPr><w:pStyle w:val="17"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr>' +
    '<w:spacing w:after="0" w:afterLines="0" w:line="240" w:lineRule="auto"/><w:ind w:left="360" w:leftChars="0"/>' +
    '<w:rPr><w:rFonts w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">- </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:b/><w:bCs w:val="0"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Artificial Intelligence:</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> Machine Learning</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
    '<w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>, Natural Language Processing</w:t></w:r>' +
    '</w:p>'
)

# 3) "- Platform: ElasticSearch, PostgreSQL, MySQL, Redis"
#    -> split out "Docker with micro-services linking, " before the rest.
Replace-ParagraphXml " ElasticSearch, PostgreSQL, MySQL, Redis" (
    '<w:p><w:pPr><w:pStyle w:val="17"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr>' +
    '<w:spacing w:after="0" w:afterLines="0" w:line="240" w:lineRule="auto"/><w:ind w:left="360" w:leftChars="0"/>' +
    '<w:rPr><w:rFonts w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">- </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:b/><w:bCs w:val="0"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Platform:</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Docker with micro-services linking, </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>ElasticSearch, PostgreSQL, MySQL, Redis</w:t></w:r>' +
    '</w:p>'
)

# 4) "- Agile Development: TDD, SCRUM"
#    -> split out "Git, " before "TDD, SCRUM".
Replace-ParagraphXml ": TDD, SCRUM" (
    '<w:p><w:pPr><w:pStyle w:val="17"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="0"/></w:numPr>' +
    '<w:spacing w:after="0" w:afterLines="0" w:line="240" w:lineRule="auto"/><w:ind w:left="360" w:leftChars="0"/>' +
    '<w:rPr><w:rFonts w:cs="Times New Roman"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">- </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:b/><w:bCs w:val="0"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Agile Development</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">Git, </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:cs="Times New Roman"/><w:b w:val="0"/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>TDD, SCRUM</w:t></w:r>' +
    '</w:p>'
)

# 5) "06/2015 - Pre[bookmark]sent Day" -> merge the split "Present Day" runs
#    back into a single run and drop the old bookmark location (it moved to #2).
Replace-ParagraphXml "Present" (
    '<w:p><w:pPr><w:wordWrap w:val="0"/><w:spacing w:after="0" w:afterLines="0" w:line="240" w:lineRule="auto"/><w:jc w:val="right"/>' +
    '<w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Microsoft Yi Baiti" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:eastAsia="Microsoft Yi Baiti" w:cs="Times New Roman"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>06</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Microsoft Yi Baiti" w:cs="Times New Roman"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>/</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:eastAsia="Microsoft Yi Baiti" w:cs="Times New Roman"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>2015</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:eastAsia="Microsoft Yi Baiti" w:cs="Times New Roman"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> - </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:eastAsia="Microsoft Yi Baiti" w:cs="Times New Roman"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Present Day</w:t></w:r>' +
    '</w:p>'
)

Write-Host "Resume skills updated (Java, Docker, Git added; Present Day run merged)."
